$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 18:04"

# Apply updated country data (values + reordering due to rank changes)
$ws.Cells.Item(4, 2).Value = 1326200  # B4: 1325791 -> 1326200
$ws.Cells.Item(4, 3).Value = 4415  # C4: 4006 -> 4415
$ws.Cells.Item(4, 5).Value = 1023471  # E4: 1023090 -> 1023471
$ws.Cells.Item(4, 7).Value = 175  # G4: 147 -> 175
$ws.Cells.Item(4, 8).Value = 78790  # H4: 78762 -> 78790
$ws.Cells.Item(35, 1).Value = "Polonia"  # A35: 'Japon' -> 'Polonia'
$ws.Cells.Item(35, 2).Value = 15651  # B35: 15575 -> 15651
$ws.Cells.Item(35, 3).Value = 285  # C35: 0 -> 285
$ws.Cells.Item(35, 4).Value = 5437  # D35: 5146 -> 5437
$ws.Cells.Item(35, 5).Value = 9429  # E35: 9839 -> 9429
$ws.Cells.Item(35, 6).Value = 160  # F35: 300 -> 160
$ws.Cells.Item(35, 7).Value = 9  # G35: 0 -> 9
$ws.Cells.Item(35, 8).Value = 785  # H35: 590 -> 785
$ws.Cells.Item(36, 1).Value = "Japon"  # A36: 'Polonia' -> 'Japon'
$ws.Cells.Item(36, 2).Value = 15575  # B36: 15510 -> 15575
$ws.Cells.Item(36, 3).Value = 0  # C36: 144 -> 0
$ws.Cells.Item(36, 4).Value = 5146  # D36: 5437 -> 5146
$ws.Cells.Item(36, 5).Value = 9839  # E36: 9296 -> 9839
$ws.Cells.Item(36, 6).Value = 300  # F36: 160 -> 300
$ws.Cells.Item(36, 7).Value = 0  # G36: 1 -> 0
$ws.Cells.Item(36, 8).Value = 590  # H36: 777 -> 590
$ws.Cells.Item(49, 2).Value = 8089  # B49: 8078 -> 8089
$ws.Cells.Item(49, 3).Value = 12  # C49: 1 -> 12
$ws.Cells.Item(49, 4).Value = 4446  # D49: 4414 -> 4446
$ws.Cells.Item(49, 5).Value = 3367  # E49: 3391 -> 3367
$ws.Cells.Item(49, 7).Value = 3  # G49: 0 -> 3
$ws.Cells.Item(49, 8).Value = 276  # H49: 273 -> 276
$ws.Cells.Item(65, 2).Value = 3877  # B65: 3871 -> 3877
$ws.Cells.Item(65, 3).Value = 6  # C65: 0 -> 6
$ws.Cells.Item(65, 4).Value = 3550  # D65: 3526 -> 3550
$ws.Cells.Item(65, 5).Value = 226  # E65: 245 -> 226
$ws.Cells.Item(65, 6).Value = 15  # F65: 16 -> 15
$ws.Cells.Item(65, 7).Value = 1  # G65: 0 -> 1
$ws.Cells.Item(65, 8).Value = 101  # H65: 100 -> 101
$ws.Cells.Item(82, 2).Value = 1754  # B82: 1741 -> 1754
$ws.Cells.Item(82, 3).Value = 13  # C82: 0 -> 13
$ws.Cells.Item(82, 4).Value = 1140  # D82: 1078 -> 1140
$ws.Cells.Item(82, 5).Value = 540  # E82: 589 -> 540
$ws.Cells.Item(121, 1).Value = "Jordania"  # A121: 'Tanzania' -> 'Jordania'
$ws.Cells.Item(121, 2).Value = 522  # B121: 509 -> 522
$ws.Cells.Item(121, 3).Value = 14  # C121: 0 -> 14
$ws.Cells.Item(121, 4).Value = 387  # D121: 183 -> 387
$ws.Cells.Item(121, 5).Value = 126  # E121: 305 -> 126
$ws.Cells.Item(121, 6).Value = 5  # F121: 7 -> 5
$ws.Cells.Item(121, 8).Value = 9  # H121: 21 -> 9
$ws.Cells.Item(122, 1).Value = "Tanzania"  # A122: 'Jordania' -> 'Tanzania'
$ws.Cells.Item(122, 2).Value = 509  # B122: 508 -> 509
$ws.Cells.Item(122, 4).Value = 183  # D122: 385 -> 183
$ws.Cells.Item(122, 5).Value = 305  # E122: 114 -> 305
$ws.Cells.Item(122, 6).Value = 7  # F122: 5 -> 7
$ws.Cells.Item(122, 8).Value = 21  # H122: 9 -> 21
$ws.Cells.Item(192, 1).Value = "Belice"  # A192: 'Nueva Caledonia' -> 'Belice'
$ws.Cells.Item(192, 4).Value = 16  # D192: 18 -> 16
$ws.Cells.Item(192, 8).Value = 2  # H192: 0 -> 2
$ws.Cells.Item(193, 1).Value = "Nueva Caledonia"  # A193: 'Belice' -> 'Nueva Caledonia'
$ws.Cells.Item(193, 4).Value = 18  # D193: 16 -> 18
$ws.Cells.Item(193, 8).Value = 0  # H193: 2 -> 0
